$wb = $excel.ActiveWorkbook

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4236.72
$ws.Cells.Item(32, 9).Value = 4171.911
$ws.Cells.Item(32, 11).Value = 4171.911
$ws.Cells.Item(32, 13).Value = -3884.911
$ws.Cells.Item(45, 8).Value = 30395980
$ws.Cells.Item(45, 9).Value = 37149976
$ws.Cells.Item(45, 10).Value = 2999
$ws.Cells.Item(45, 11).Value = 37149976
$ws.Cells.Item(45, 12).Value = 2999
$ws.Cells.Item(45, 13).Value = -37149599
$ws.Cells.Item(45, 14).Value = -3753
$ws.Cells.Item(61, 8).Value = 1815.7333
$ws.Cells.Item(61, 9).Value = 1853.2174
$ws.Cells.Item(61, 10).Value = 1776.5454
$ws.Cells.Item(61, 11).Value = 1853.2174
$ws.Cells.Item(61, 12).Value = 1776.5454
$ws.Cells.Item(61, 13).Value = -1641.2174
$ws.Cells.Item(61, 14).Value = -2200.5454
$ws.Cells.Item(110, 8).Value = 1591.091
$ws.Cells.Item(110, 9).Value = 936.7059
$ws.Cells.Item(110, 10).Value = 3816
$ws.Cells.Item(110, 11).Value = 936.7059
$ws.Cells.Item(110, 12).Value = 3816
$ws.Cells.Item(110, 13).Value = 1108.2941
$ws.Cells.Item(110, 14).Value = -7906
$ws.Cells.Item(122, 8).Value = 1299.9
$ws.Cells.Item(122, 9).Value = 1222.1111
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 3666.3333
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -1216.3333
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(126, 8).Value = 5680
$ws.Cells.Item(126, 9).Value = 5680
$ws.Cells.Item(126, 11).Value = 17040
$ws.Cells.Item(126, 13).Value = -14570
$ws.Cells.Item(132, 8).Value = 23813310
$ws.Cells.Item(132, 9).Value = 30304170
$ws.Cells.Item(132, 10).Value = 13487.889
$ws.Cells.Item(132, 11).Value = 90912510
$ws.Cells.Item(132, 12).Value = 40463.667
$ws.Cells.Item(132, 13).Value = -90909980
$ws.Cells.Item(132, 14).Value = -45523.667
$ws.Cells.Item(136, 8).Value = 1815.7333
$ws.Cells.Item(136, 9).Value = 1853.2174
$ws.Cells.Item(136, 10).Value = 1776.5454
$ws.Cells.Item(136, 11).Value = 5559.6522
$ws.Cells.Item(136, 12).Value = 5329.6362
$ws.Cells.Item(136, 13).Value = -3009.6522
$ws.Cells.Item(136, 14).Value = -10429.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2140.647
$ws.Cells.Item(99, 9).Value = 1050
$ws.Cells.Item(99, 10).Value = 2286.0667
$ws.Cells.Item(99, 11).Value = 1050
$ws.Cells.Item(99, 12).Value = 2286.0667
$ws.Cells.Item(99, 13).Value = 448
$ws.Cells.Item(99, 14).Value = -5282.066699999999
$ws.Cells.Item(107, 8).Value = 10870256
$ws.Cells.Item(107, 9).Value = 12821084
$ws.Cells.Item(107, 11).Value = 12821084
$ws.Cells.Item(107, 13).Value = -12819164
$ws.Cells.Item(134, 8).Value = 3007612.2
$ws.Cells.Item(134, 9).Value = 1349.762
$ws.Cells.Item(134, 11).Value = 4049.286
$ws.Cells.Item(134, 13).Value = -1514.286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1480.34
$ws.Cells.Item(31, 9).Value = 1011.4286
$ws.Cells.Item(31, 10).Value = 2278.7568
$ws.Cells.Item(31, 11).Value = 1011.4286
$ws.Cells.Item(31, 12).Value = 2278.7568
$ws.Cells.Item(31, 13).Value = -716.4286
$ws.Cells.Item(31, 14).Value = -2868.7568
$ws.Cells.Item(34, 8).Value = 1480.34
$ws.Cells.Item(34, 9).Value = 1011.4286
$ws.Cells.Item(34, 10).Value = 2278.7568
$ws.Cells.Item(34, 11).Value = 1011.4286
$ws.Cells.Item(34, 12).Value = 2278.7568
$ws.Cells.Item(34, 13).Value = -809.4286
$ws.Cells.Item(34, 14).Value = -2682.7568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 922.4299999999999
$ws.Cells.Item(131, 10).Value = 924.6701
$ws.Cells.Item(131, 12).Value = 2774.0103
$ws.Cells.Item(131, 14).Value = -12854.0103

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(6, 8).Value = 11000
$ws.Cells.Item(6, 10).Value = 11000
$ws.Cells.Item(6, 12).Value = 11000
$ws.Cells.Item(6, 14).Value = -11226
$ws.Cells.Item(7, 8).Value = 5000000
$ws.Cells.Item(7, 9).Value = 5000000
$ws.Cells.Item(7, 10).Value = 5000000
$ws.Cells.Item(7, 11).Value = 5000000
$ws.Cells.Item(7, 12).Value = 5000000
$ws.Cells.Item(7, 13).Value = -4999888
$ws.Cells.Item(7, 14).Value = -5000224
$ws.Cells.Item(8, 8).Value = 5000000
$ws.Cells.Item(8, 9).Value = 5000000
$ws.Cells.Item(8, 10).Value = 5000000
$ws.Cells.Item(8, 11).Value = 5000000
$ws.Cells.Item(8, 12).Value = 5000000
$ws.Cells.Item(8, 13).Value = -4999861
$ws.Cells.Item(8, 14).Value = -5000278
$ws.Cells.Item(10, 8).Value = 5000434.5
$ws.Cells.Item(10, 9).Value = 7500251.5
$ws.Cells.Item(10, 10).Value = 800
$ws.Cells.Item(10, 11).Value = 7500251.5
$ws.Cells.Item(10, 12).Value = 800
$ws.Cells.Item(10, 13).Value = -7500082.5
$ws.Cells.Item(10, 14).Value = -1138
$ws.Cells.Item(11, 8).Value = 3498796.2
$ws.Cells.Item(11, 9).Value = 4081766.8
$ws.Cells.Item(11, 10).Value = 974
$ws.Cells.Item(11, 11).Value = 4081766.8
$ws.Cells.Item(11, 12).Value = 974
$ws.Cells.Item(11, 13).Value = -4081627.8
$ws.Cells.Item(11, 14).Value = -1252
$ws.Cells.Item(13, 8).Value = 570
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(14, 8).Value = 250286.67
$ws.Cells.Item(14, 9).Value = 333537.78
$ws.Cells.Item(14, 10).Value = 533.3333
$ws.Cells.Item(14, 11).Value = 333537.78
$ws.Cells.Item(14, 12).Value = 533.3333
$ws.Cells.Item(14, 13).Value = -333369.78
$ws.Cells.Item(14, 14).Value = -869.3333
$ws.Cells.Item(16, 8).Value = 11000
$ws.Cells.Item(16, 10).Value = 11000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 14).Value = -11500
$ws.Cells.Item(17, 8).Value = 3000
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 3000
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 3000
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).Value = -3336
$ws.Cells.Item(19, 8).Value = 987.6667
$ws.Cells.Item(19, 9).Value = 987.6667
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 987.6667
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -699.6667
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(21, 8).Value = 24980
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()
$ws.Cells.Item(24, 8).Value = 100
$ws.Cells.Item(24, 9).Value = 100
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 100
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 73
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 5000
$ws.Cells.Item(27, 10).Value = 5000
$ws.Cells.Item(27, 12).Value = 5000
$ws.Cells.Item(27, 14).Value = -5332
$ws.Cells.Item(28, 8).Value = 15000
$ws.Cells.Item(28, 10).Value = 15000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 14).Value = -15384
$ws.Cells.Item(30, 8).Value = 24980
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 13).ClearContents()
$ws.Cells.Item(31, 8).Value = 931
$ws.Cells.Item(31, 9).Value = 931
$ws.Cells.Item(31, 11).Value = 931
$ws.Cells.Item(31, 13).Value = -639
$ws.Cells.Item(36, 8).Value = 776
$ws.Cells.Item(36, 9).Value = 659
$ws.Cells.Item(36, 10).Value = 1010
$ws.Cells.Item(36, 11).Value = 659
$ws.Cells.Item(36, 12).Value = 1010
$ws.Cells.Item(36, 13).Value = -174
$ws.Cells.Item(36, 14).Value = -1980
$ws.Cells.Item(37, 8).Value = 931
$ws.Cells.Item(37, 9).Value = 931
$ws.Cells.Item(37, 11).Value = 931
$ws.Cells.Item(37, 13).Value = -654
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(41, 8).Value = 10350.333
$ws.Cells.Item(41, 9).Value = 675.5
$ws.Cells.Item(41, 10).Value = 29700
$ws.Cells.Item(41, 11).Value = 675.5
$ws.Cells.Item(41, 12).Value = 29700
$ws.Cells.Item(41, 13).Value = -320.5
$ws.Cells.Item(41, 14).Value = -30410
$ws.Cells.Item(43, 8).Value = 3579.75
$ws.Cells.Item(43, 10).Value = 8019
$ws.Cells.Item(43, 12).Value = 8019
$ws.Cells.Item(43, 14).Value = -8321
$ws.Cells.Item(48, 8).Value = 5000
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(49, 8).Value = 6000
$ws.Cells.Item(49, 10).Value = 6000
$ws.Cells.Item(49, 12).Value = 6000
$ws.Cells.Item(49, 14).Value = -6368
$ws.Cells.Item(53, 8).Value = 4125
$ws.Cells.Item(53, 9).Value = 2000
$ws.Cells.Item(53, 10).Value = 4833.3335
$ws.Cells.Item(53, 11).Value = 2000
$ws.Cells.Item(53, 12).Value = 4833.3335
$ws.Cells.Item(53, 13).Value = -1369
$ws.Cells.Item(53, 14).Value = -6095.3335
$ws.Cells.Item(55, 8).Value = 9750
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 9750
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 9750
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -10404
$ws.Cells.Item(57, 8).Value = 14020.333
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(58, 8).Value = 166671170
$ws.Cells.Item(58, 9).Value = 2000
$ws.Cells.Item(58, 10).Value = 200005000
$ws.Cells.Item(58, 11).Value = 2000
$ws.Cells.Item(58, 12).Value = 200005000
$ws.Cells.Item(58, 13).Value = -1723
$ws.Cells.Item(58, 14).Value = -200005554
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 33343650
$ws.Cells.Item(122, 9).Value = 45466430
$ws.Cells.Item(122, 10).Value = 6000.25
$ws.Cells.Item(122, 11).Value = 136399290
$ws.Cells.Item(122, 12).Value = 18000.75
$ws.Cells.Item(122, 13).Value = -136396840
$ws.Cells.Item(122, 14).Value = -22900.75
$ws.Cells.Item(132, 8).Value = 6315.8184
$ws.Cells.Item(132, 9).Value = 1356.6666
$ws.Cells.Item(132, 11).Value = 4069.9998
$ws.Cells.Item(132, 13).Value = -1539.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 50002420
$ws.Cells.Item(40, 9).Value = 3026
$ws.Cells.Item(40, 11).Value = 3026
$ws.Cells.Item(40, 13).Value = -2890

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 12532.4
$ws.Cells.Item(64, 10).Value = 12532.4
$ws.Cells.Item(64, 12).Value = 12532.4
$ws.Cells.Item(64, 14).Value = -13028.4
$ws.Cells.Item(67, 8).Value = 12532.4
$ws.Cells.Item(67, 10).Value = 12532.4
$ws.Cells.Item(67, 12).Value = 12532.4
$ws.Cells.Item(67, 14).Value = -14248.4
$ws.Cells.Item(132, 8).Value = 26677.63
$ws.Cells.Item(132, 9).Value = 38578.18
$ws.Cells.Item(132, 10).Value = 8165.6665
$ws.Cells.Item(132, 11).Value = 115734.54
$ws.Cells.Item(132, 12).Value = 24496.9995
$ws.Cells.Item(132, 13).Value = -113204.54
$ws.Cells.Item(132, 14).Value = -29556.9995
